$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix header labels: column C was "sum" -> "nombre_aides"; column D was "nombre_aides" -> "montant_total"
$ws.Range("C1").Value2 = "nombre_aides"
$ws.Range("D1").Value2 = "montant_total"

# Find the last used row
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# Swap the values of columns C and D for every data row (C held the monetary sum, D held the count;
# now C should hold the count and D should hold the monetary sum)
for ($r = 2; $r -le $lastRow; $r++) {
    $cVal = $ws.Cells.Item($r, 3).Value2
    $dVal = $ws.Cells.Item($r, 4).Value2
    $ws.Cells.Item($r, 3).Value2 = $dVal
    $ws.Cells.Item($r, 4).Value2 = $cVal
}

# Fix the mojibake text in column A ("Fonds de solidarit??" -> "Fonds de solidarité")
for ($r = 2; $r -le $lastRow; $r++) {
    $aVal = $ws.Cells.Item($r, 1).Value2
    if ($aVal -eq "Fonds de solidarit??") {
        $ws.Cells.Item($r, 1).Value2 = "Fonds de solidarité"
    }
}
